$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.34
$ws.Range("K2").Value = 4.1
$ws.Range("N2").Value = 4.7
$ws.Range("O2").Value = 1.24
$ws.Range("P2").Value = 2.26
$ws.Range("Q2").Value = 1.72
$ws.Range("R2").Value = 1.51
$ws.Range("S2").Value = 2.8
$ws.Range("V2").Value = 1.43
$ws.Range("W2").Value = 1.75
$ws.Range("AK2").Value = 27
$ws.Range("AL2").Value = 38
$ws.Range("AN2").Value = 16.5
$ws.Range("Q3").Value = 1.87
$ws.Range("F4").Value = 2.18
$ws.Range("G4").Value = 2.5
$ws.Range("H4").Value = 3.55
$ws.Range("I4").Value = 4.6
$ws.Range("J4").Value = 3.15
$ws.Range("K4").Value = 3.9
$ws.Range("P4").Value = 1.65
$ws.Range("F5").Value = 3.6
$ws.Range("G5").Value = 5.4
$ws.Range("H5").Value = 1.82
$ws.Range("I5").Value = 2.26
$ws.Range("J5").Value = 3.2
$ws.Range("K5").Value = 5.9
$ws.Range("N5").Value = 1.86
$ws.Range("P5").Value = 1.86
$ws.Range("Q5").Value = 1.74
$ws.Range("S5").Value = 2.78
$ws.Range("V5").Value = 1.79
$ws.Range("W5").Value = 1.23
$ws.Range("G6").Value = 9.199999999999999
$ws.Range("K6").Value = 5.4
$ws.Range("Q6").Value = 1.69
$ws.Range("P7").Value = 3.15
$ws.Range("F8").Value = 1.33
$ws.Range("G8").Value = 1.41
$ws.Range("H8").Value = 8.199999999999999
$ws.Range("I8").Value = 10.5
$ws.Range("J8").Value = 6
$ws.Range("K8").Value = 6.8
$ws.Range("O8").Value = 1.14
$ws.Range("P8").Value = 2.94
$ws.Range("R8").Value = 1.79
$ws.Range("S8").Value = 2.02
$ws.Range("Z8").Value = 110
$ws.Range("AA8").Value = 320
$ws.Range("AB8").Value = 14.5
$ws.Range("AF8").Value = 11.5
$ws.Range("AI8").Value = 100
$ws.Range("F10").Value = 3.3
$ws.Range("H10").Value = 2.08
$ws.Range("I10").Value = 2.66
$ws.Range("J10").Value = 2.84
$ws.Range("K10").Value = 4.7
$ws.Range("P10").Value = 1.69
$ws.Range("Q10").Value = 1.91
$ws.Range("F11").Value = 1.81
$ws.Range("G11").Value = 2.26
$ws.Range("I11").Value = 4.8
$ws.Range("J11").Value = 2.86
$ws.Range("K11").Value = 7.2
$ws.Range("I12").Value = 3.9
$ws.Range("J12").Value = 3.8
$ws.Range("Q12").Value = 1.37
$ws.Range("G14").Value = 1.99
$ws.Range("H14").Value = 3.2
$ws.Range("J14").Value = 3.5
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 3.35
$ws.Range("H15").Value = 2.12
$ws.Range("I15").Value = 2.3
$ws.Range("K15").Value = 4.6
$ws.Range("P15").Value = 2.82
$ws.Range("F20").Value = 1.23
$ws.Range("H20").Value = 3.55
$ws.Range("K20").Value = 980
$ws.Range("Q20").Value = 1.01
$ws.Range("F21").Value = 2.6
$ws.Range("G21").Value = 2.92
$ws.Range("I21").Value = 3.55
$ws.Range("J21").Value = 2.96
$ws.Range("AF22").Value = 22
$ws.Range("R24").Value = 1.27
$ws.Range("AJ24").Value = 40
$ws.Range("AO24").Value = 46
$ws.Range("AE25").Value = 1000
$ws.Range("I26").Value = 14.5
$ws.Range("P26").Value = 1.82
